$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "User Name"
$ws.Range("B1").Value = "Computer Name"
$ws.Range("C1").Value = "Static IP"

# Data rows
$data = @(
    @("jacky", "test4", "192.168.11.23"),
    @("jacky", "test41", "192.168.11.231"),
    @("jacky", "test42", "192.168.11.232"),
    @("jacky", "test", "192.168.11.22")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Header styling: bold font, centered horizontal, top vertical align, thin border all around
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous
$headerRange.Borders.Weight = 2           # xlThin
